$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 32; $row++) {
    $formula = '=CONCATENATE("EMOTION(",H' + $row + ',"){value=[",B' + $row + ',"]; value_type=[BASICEMOTION]; emotionIntensity=[",SUBSTITUTE(C' + $row + ',",","."),"]; sourceAggr=[",SUBSTITUTE(D' + $row + ',",","."),"]; sourceLibid=[",SUBSTITUTE(E' + $row + ',",","."),"]; sourcePleasure=[",SUBSTITUTE(F' + $row + ',",","."),"]; sourceUnpleasure=[",SUBSTITUTE(G' + $row + ',",","."),"]}")'
    $ws.Range("I$row").Formula = $formula
}

$ws.Range("J2:J32").Select()
